$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update row 2 values (B2:E2)
$ws.Range("B2").Value = 4.4286150598072966
$ws.Range("C2").Value = 4.740854810894489
$ws.Range("D2").Value = 6.5584700280543338
$ws.Range("E2").Value = 5.4368871074205849

# Update row 3 values (B3:E3)
$ws.Range("B3").Value = 4.8852490163363234
$ws.Range("C3").Value = 6.0333219377177141
$ws.Range("D3").Value = 5.5878078263626962
$ws.Range("E3").Value = 5.202518717321241

# Update selection to match new used focus range
$ws.Range("B1:E3").Select()
